# Update "Pais" COVID data sheet: refresh the "last updated" timestamp and
# apply the latest per-country counters. Some countries changed rank (their
# "Casos totales" overtook/fell behind a neighboring country), so those rows
# swap their country name + figures with the adjacent row to keep the table
# sorted by "Casos totales" descending, exactly like the source refresh did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Junio de 2020 a las 06:07"

# Row 7: India - update Casos activos / Recuperados
$ws.Cells.Item(7,1).Value = "India"
$ws.Cells.Item(7,2).Value = 321626
$ws.Cells.Item(7,3).Value = 0
$ws.Cells.Item(7,4).Value = 162379
$ws.Cells.Item(7,5).Value = 150048
$ws.Cells.Item(7,6).Value = 0
$ws.Cells.Item(7,7).Value = 0
$ws.Cells.Item(7,8).Value = 9199

# Row 11: Peru - update Casos totales / activos / recuperados / muertes
$ws.Cells.Item(11,1).Value = "Peru"
$ws.Cells.Item(11,2).Value = 225132
$ws.Cells.Item(11,3).Value = 0
$ws.Cells.Item(11,4).Value = 111724
$ws.Cells.Item(11,5).Value = 106910
$ws.Cells.Item(11,6).Value = 0
$ws.Cells.Item(11,7).Value = 0
$ws.Cells.Item(11,8).Value = 6498

# Row 56: Kazajistan - update totales / nuevos / activos / recuperados
$ws.Cells.Item(56,1).Value = "Kazajistan"
$ws.Cells.Item(56,2).Value = 14496
$ws.Cells.Item(56,3).Value = 258
$ws.Cells.Item(56,4).Value = 9056
$ws.Cells.Item(56,5).Value = 5367
$ws.Cells.Item(56,6).Value = 0
$ws.Cells.Item(56,7).Value = 0
$ws.Cells.Item(56,8).Value = 73

# Rows 69-70: Honduras overtakes Malasia (Honduras rises above, Malasia keeps its
# previous figures but drops one row)
$ws.Cells.Item(69,1).Value = "Honduras"
$ws.Cells.Item(69,2).Value = 8455
$ws.Cells.Item(69,3).Value = 323
$ws.Cells.Item(69,4).Value = 894
$ws.Cells.Item(69,5).Value = 7251
$ws.Cells.Item(69,6).Value = 0
$ws.Cells.Item(69,7).Value = 4
$ws.Cells.Item(69,8).Value = 310

$ws.Cells.Item(70,1).Value = "Malasia"
$ws.Cells.Item(70,2).Value = 8445
$ws.Cells.Item(70,3).Value = 0
$ws.Cells.Item(70,4).Value = 7311
$ws.Cells.Item(70,5).Value = 1014
$ws.Cells.Item(70,6).Value = 0
$ws.Cells.Item(70,7).Value = 0
$ws.Cells.Item(70,8).Value = 120

# Rows 96-99: Kirguistan overtakes Mayotte, Croacia and Cuba, which each drop
# one row keeping their previous figures
$ws.Cells.Item(96,1).Value = "Kirguistan"
$ws.Cells.Item(96,2).Value = 2285
$ws.Cells.Item(96,3).Value = 78
$ws.Cells.Item(96,4).Value = 1791
$ws.Cells.Item(96,5).Value = 467
$ws.Cells.Item(96,6).Value = 0
$ws.Cells.Item(96,7).Value = 0
$ws.Cells.Item(96,8).Value = 27

$ws.Cells.Item(97,1).Value = "Mayotte"
$ws.Cells.Item(97,2).Value = 2282
$ws.Cells.Item(97,3).Value = 0
$ws.Cells.Item(97,4).Value = 1790
$ws.Cells.Item(97,5).Value = 464
$ws.Cells.Item(97,6).Value = 0
$ws.Cells.Item(97,7).Value = 0
$ws.Cells.Item(97,8).Value = 28

$ws.Cells.Item(98,1).Value = "Croacia"
$ws.Cells.Item(98,2).Value = 2251
$ws.Cells.Item(98,3).Value = 0
$ws.Cells.Item(98,4).Value = 2134
$ws.Cells.Item(98,5).Value = 10
$ws.Cells.Item(98,6).Value = 0
$ws.Cells.Item(98,7).Value = 0
$ws.Cells.Item(98,8).Value = 107

$ws.Cells.Item(99,1).Value = "Cuba"
$ws.Cells.Item(99,2).Value = 2238
$ws.Cells.Item(99,3).Value = 0
$ws.Cells.Item(99,4).Value = 1923
$ws.Cells.Item(99,5).Value = 231
$ws.Cells.Item(99,6).Value = 0
$ws.Cells.Item(99,7).Value = 0
$ws.Cells.Item(99,8).Value = 84

# Row 162: Mongolia - update Casos activos / Recuperados
$ws.Cells.Item(162,1).Value = "Mongolia"
$ws.Cells.Item(162,2).Value = 197
$ws.Cells.Item(162,3).Value = 0
$ws.Cells.Item(162,4).Value = 98
$ws.Cells.Item(162,5).Value = 99
$ws.Cells.Item(162,6).Value = 0
$ws.Cells.Item(162,7).Value = 0
$ws.Cells.Item(162,8).Value = 0

# Rows 183-184: Butan overtakes Eritrea (Butan rises above, Eritrea keeps its
# previous figures but drops one row)
$ws.Cells.Item(183,1).Value = "Butan"
$ws.Cells.Item(183,2).Value = 66
$ws.Cells.Item(183,3).Value = 4
$ws.Cells.Item(183,4).Value = 21
$ws.Cells.Item(183,5).Value = 45
$ws.Cells.Item(183,6).Value = 0
$ws.Cells.Item(183,7).Value = 0
$ws.Cells.Item(183,8).Value = 0

$ws.Cells.Item(184,1).Value = "Eritrea"
$ws.Cells.Item(184,2).Value = 65
$ws.Cells.Item(184,3).Value = 0
$ws.Cells.Item(184,4).Value = 39
$ws.Cells.Item(184,5).Value = 26
$ws.Cells.Item(184,6).Value = 0
$ws.Cells.Item(184,7).Value = 0
$ws.Cells.Item(184,8).Value = 0

# Rows 213-214: Papua Nueva Guinea and Islas Virgenes Britanicas swap order
# (tied Casos totales, figures unchanged, only row order flips)
$ws.Cells.Item(213,1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213,2).Value = 8
$ws.Cells.Item(213,3).Value = 0
$ws.Cells.Item(213,4).Value = 8
$ws.Cells.Item(213,5).Value = 0
$ws.Cells.Item(213,6).Value = 0
$ws.Cells.Item(213,7).Value = 0
$ws.Cells.Item(213,8).Value = 0

$ws.Cells.Item(214,1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214,2).Value = 8
$ws.Cells.Item(214,3).Value = 0
$ws.Cells.Item(214,4).Value = 7
$ws.Cells.Item(214,5).Value = 0
$ws.Cells.Item(214,6).Value = 0
$ws.Cells.Item(214,7).Value = 0
$ws.Cells.Item(214,8).Value = 1
